# Auto-generated Excel COM-interop script
# Applies the value changes from the commit diff to Sheets/Hyperion_Profits.xlsx
# (workbook sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 3450
$ws.Range("I32").Value = 3450
$ws.Range("K32").Value = 3450
$ws.Range("M32").Value = -3124
# Row 43
$ws.Range("H43").Value = 4357.7144
$ws.Range("I43").Value = 1813.5
$ws.Range("J43").Value = 7750
$ws.Range("K43").Value = 1813.5
$ws.Range("L43").Value = 7750
$ws.Range("M43").Value = -1744.5
$ws.Range("N43").Value = -7888
# Row 54
$ws.Range("H54").Value = 5076
$ws.Range("I54").Value = 5076
$ws.Range("K54").Value = 5076
$ws.Range("M54").Value = -4590
# Row 112
$ws.Range("H112").Value = 5253.846
$ws.Range("J112").Value = 5430
$ws.Range("L112").Value = 16290
$ws.Range("N112").Value = -18506
# Row 115
$ws.Range("H115").Value = 1011.53845
$ws.Range("I115").Value = 1011.53845
$ws.Range("K115").Value = 3034.61535
$ws.Range("M115").Value = -1467.61535
# Row 116
$ws.Range("H116").Value = 4578.5
$ws.Range("I116").Value = 4178.467
$ws.Range("K116").Value = 4178.467
$ws.Range("M116").Value = -736.4669999999996
# Row 132
$ws.Range("H132").Value = 31253992
$ws.Range("I132").Value = 33337532
$ws.Range("K132").Value = 100012596
$ws.Range("M132").Value = -100010066

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 14
$ws.Range("H14").Value = 601.2
$ws.Range("I14").Value = 468.66666
$ws.Range("J14").Value = 800
$ws.Range("K14").Value = 468.66666
$ws.Range("L14").Value = 800
$ws.Range("M14").Value = -293.66666
$ws.Range("N14").Value = -1150
# Row 45
$ws.Range("H45").Value = 7196484
$ws.Range("I45").Value = 11989717
$ws.Range("J45").Value = 6634.5
$ws.Range("K45").Value = 11989717
$ws.Range("L45").Value = 6634.5
$ws.Range("M45").Value = -11989340
$ws.Range("N45").Value = -7388.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 1067.6666
$ws.Range("I7").Value = 1067.6666
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1067.6666
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -954.6666
$ws.Range("N7").ClearContents()
# Row 97
$ws.Range("H97").Value = 15375.333
$ws.Range("I97").Value = 15375.333
$ws.Range("K97").Value = 15375.333
$ws.Range("M97").Value = -14384.333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 4705.643
$ws.Range("I134").Value = 3236.125
$ws.Range("J134").Value = 6665
$ws.Range("K134").Value = 9708.375
$ws.Range("L134").Value = 19995
$ws.Range("M134").Value = -7173.375
$ws.Range("N134").Value = -25065
# Row 141
$ws.Range("H141").Value = 408289.8
$ws.Range("J141").Value = 408289.8
$ws.Range("L141").Value = 408289.8
$ws.Range("N141").Value = -418649.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 103
$ws.Range("H103").Value = 366
$ws.Range("I103").Value = 239.2
$ws.Range("K103").Value = 717.5999999999999
$ws.Range("M103").Value = 161.4000000000001
# Row 115
$ws.Range("H115").Value = 1500
$ws.Range("I115").Value = 1500
$ws.Range("K115").Value = 4500
$ws.Range("M115").Value = -3325
# Row 117
$ws.Range("H117").Value = 835.6667
$ws.Range("I117").Value = 160.66667
$ws.Range("J117").Value = 1004.4167
$ws.Range("K117").Value = 482.00001
$ws.Range("L117").Value = 3013.2501
$ws.Range("M117").Value = 2959.99999
$ws.Range("N117").Value = -9897.250100000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 11285.643
$ws.Range("I43").Value = 5124.875
$ws.Range("K43").Value = 5124.875
$ws.Range("M43").Value = -4973.875
# Row 46
$ws.Range("H46").Value = 6960.6665
$ws.Range("I46").Value = 6960.6665
$ws.Range("K46").Value = 6960.6665
$ws.Range("M46").Value = -6804.6665
# Row 57
$ws.Range("H57").Value = 14068.182
$ws.Range("J57").Value = 14068.182
$ws.Range("L57").Value = 14068.182
$ws.Range("N57").Value = -15708.182
# Row 80
$ws.Range("H80").Value = 27606798
$ws.Range("J80").Value = 2615.6667
$ws.Range("L80").Value = 2615.6667
$ws.Range("N80").Value = -4611.6667
# Row 83
$ws.Range("H83").Value = 27606798
$ws.Range("J83").Value = 2615.6667
$ws.Range("L83").Value = 13078.3335
$ws.Range("N83").Value = -23062.3335
# Row 141
$ws.Range("H141").Value = 54500
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 54500
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 54500
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -64860

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 9120.143
$ws.Range("I7").Value = 5806.4
$ws.Range("J7").Value = 10961.111
$ws.Range("K7").Value = 5806.4
$ws.Range("L7").Value = 10961.111
$ws.Range("M7").Value = -5694.4
$ws.Range("N7").Value = -11185.111
# Row 19
$ws.Range("H19").Value = 1830
$ws.Range("I19").Value = 399.33334
$ws.Range("J19").Value = 2688.4
$ws.Range("K19").Value = 399.33334
$ws.Range("L19").Value = 2688.4
$ws.Range("M19").Value = -229.33334
$ws.Range("N19").Value = -3028.4
# Row 22
$ws.Range("H22").Value = 76148.336
$ws.Range("I22").Value = 149646.33
$ws.Range("J22").Value = 2650.3333
$ws.Range("K22").Value = 149646.33
$ws.Range("L22").Value = 2650.3333
$ws.Range("M22").Value = -149351.33
$ws.Range("N22").Value = -3240.3333
# Row 27
$ws.Range("H27").Value = 76148.336
$ws.Range("I27").Value = 149646.33
$ws.Range("J27").Value = 2650.3333
$ws.Range("K27").Value = 149646.33
$ws.Range("L27").Value = 2650.3333
$ws.Range("M27").Value = -149539.33
$ws.Range("N27").Value = -2864.3333
# Row 40
$ws.Range("H40").Value = 4455.3955
$ws.Range("I40").Value = 3500.3333
$ws.Range("K40").Value = 3500.3333
$ws.Range("M40").Value = -3364.3333
# Row 62
$ws.Range("H62").Value = 55530.5
$ws.Range("J62").Value = 55530.5
$ws.Range("L62").Value = 55530.5
$ws.Range("N62").Value = -56778.5
# Row 65
$ws.Range("H65").Value = 55530.5
$ws.Range("J65").Value = 55530.5
$ws.Range("L65").Value = 166591.5
$ws.Range("N65").Value = -172831.5
# Row 94
$ws.Range("H94").Value = 46500
$ws.Range("J94").Value = 46500
$ws.Range("L94").Value = 46500
$ws.Range("N94").Value = -47852
# Row 126
$ws.Range("H126").Value = 9120.143
$ws.Range("I126").Value = 5806.4
$ws.Range("J126").Value = 10961.111
$ws.Range("K126").Value = 17419.2
$ws.Range("L126").Value = 32883.333
$ws.Range("M126").Value = -14949.2
$ws.Range("N126").Value = -37823.333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 3333
$ws.Range("J2").Value = 3333
$ws.Range("L2").Value = 3333
$ws.Range("N2").Value = -3557
# Row 3
$ws.Range("H3").Value = 2386
$ws.Range("J3").Value = 2386
$ws.Range("L3").Value = 2386
$ws.Range("N3").Value = -2614
# Row 4
$ws.Range("H4").Value = 313481.34
$ws.Range("I4").Value = 467500
$ws.Range("J4").Value = 5444
$ws.Range("K4").Value = 467500
$ws.Range("L4").Value = 5444
$ws.Range("M4").Value = -467387
$ws.Range("N4").Value = -5670
# Row 10
$ws.Range("H10").Value = 1829
$ws.Range("I10").Value = 325
$ws.Range("J10").Value = 3333
$ws.Range("K10").Value = 325
$ws.Range("L10").Value = 3333
$ws.Range("M10").Value = -156
$ws.Range("N10").Value = -3671
# Row 14
$ws.Range("H14").Value = 10005
$ws.Range("J14").Value = 10005
$ws.Range("L14").Value = 10005
$ws.Range("N14").Value = -10341
# Row 122
$ws.Range("H122").Value = 1821.4222
$ws.Range("I122").Value = 1650.6923
$ws.Range("J122").Value = 2055.0527
$ws.Range("K122").Value = 4952.0769
$ws.Range("L122").Value = 6165.158100000001
$ws.Range("M122").Value = -2502.0769
$ws.Range("N122").Value = -11065.1581
# Row 136
$ws.Range("H136").Value = 4312.6113
$ws.Range("J136").Value = 5259.385
$ws.Range("L136").Value = 15778.155
$ws.Range("N136").Value = -20878.155
